# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right before the "总计" (total) sheet,
#    with the same per-quarter fund-holding layout used by the other
#    quarter sheets, populated with the new holding snapshot.
# 2. Rebuild the "总计" (total) sheet with a new leading row for 2022-Q1
#    and every other row shifted down by one (counter column updated to
#    match), re-using the formatting of the surviving rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as genuine TEXT (not auto-coerced to
# a number), without leaving a stray NumberFormat/style behind.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: reorder/insert sheets so the final tab order + sheetId layout
# is 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
# (mirrors how Excel assigns a fresh sheetId to every newly Add()-ed
# sheet, so the old "总计" is dropped and re-created last).
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

$q1sheet = $wb.Worksheets.Add($null, $template)
$q1sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add($null, $q1sheet)
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: populate the new "2022-Q1" sheet (same shape as the other
# quarter sheets: fund code / name / size / total equity position /
# position weight / held market value / position rank).
# ---------------------------------------------------------------------
$template.Range("B1:H1").Copy()
$q1sheet.Range("B1").PasteSpecial(-4122) # xlPasteFormats
$template.Range("A2").Copy()
$q1sheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$q1sheet.Range("B1").Value = "基金代码"
$q1sheet.Range("C1").Value = "基金名称"
$q1sheet.Range("D1").Value = "基金规模"
$q1sheet.Range("E1").Value = "股票总仓位"
$q1sheet.Range("F1").Value = "仓位占比"
$q1sheet.Range("G1").Value = "持有市值(亿元)"
$q1sheet.Range("H1").Value = "仓位排名"

$q1sheet.Range("A2").Value = 0
Set-TextValue $q1sheet.Range("B2") "513080"
$q1sheet.Range("C2").Value = "华安法国CAC40ETF（QDII）"
Set-TextValue $q1sheet.Range("D2") "0.60"
Set-TextValue $q1sheet.Range("E2") "96.69"
Set-TextValue $q1sheet.Range("F2") "5.64"
Set-TextValue $q1sheet.Range("G2") "0.0338"
$q1sheet.Range("H2").Value = 4

# ---------------------------------------------------------------------
# Step 3: rebuild the "总计" (total) sheet: header + one row per quarter,
# newest first, with a 0-based running index in column A.
# ---------------------------------------------------------------------
$template.Range("B1:D1").Copy()
$totalSheet.Range("B1").PasteSpecial(-4122) # xlPasteFormats
$template.Range("A2").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122) # xlPasteFormats

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$quarters = @(
    @{ Label = "2022-Q1"; Count = 1; Value = 0.03 },
    @{ Label = "2021-Q4"; Count = 1; Value = 0.04 },
    @{ Label = "2021-Q3"; Count = 1; Value = 0.03 },
    @{ Label = "2021-Q2"; Count = 1; Value = 0.03 },
    @{ Label = "2021-Q1"; Count = 1; Value = 0.03 },
    @{ Label = "2020-Q4"; Count = 1; Value = 0.03 }
)

for ($i = 0; $i -lt $quarters.Count; $i++) {
    $row = $i + 2
    $totalSheet.Range("A$row").Value = $i
    $totalSheet.Range("B$row").Value = $quarters[$i].Label
    $totalSheet.Range("C$row").Value = $quarters[$i].Count
    $totalSheet.Range("D$row").Value = $quarters[$i].Value
}

# Restore the original active sheet/selection (deleting + re-adding
# "总计" shifts the active tab, which isn't otherwise part of this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
